# Generate Report for Handoff
#
# The localization pipeline moved this file from "In Translation" to
# "Ready for handoff": update the status text on every sheet that shows
# it, bump the two "generate/handoff datetime" stamps forward a minute,
# and let the (now-wider) Status-ish columns re-size to fit the new text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# Columns E/F are the per-language status ("zh-cn" / "de-de"); G is the
# latest handoff-xliff generation timestamp.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-14 02:54:53"

# --- zh-cn detail sheet ------------------------------------------------
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-14 02:54:44"

# --- de-de detail sheet ------------------------------------------------
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-14 02:54:53"

# --- Resize the status columns for the new, longer text ---------------
# "Ready for handoff" is wider than "In Translation", so Excel widens
# the holding column when the report is regenerated.
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
